# Ordered (old, new) text pairs for each of the 100 table cells (row-major),
# matching the document's cell order exactly (includes two cells that both
# originally read "91-85=6" but resolve to different replacements).
$pairs = @(
    @("84+3=87", "73-45=28"),
    @("40+26=66", "78-63=15"),
    @("82-77=5", "21+23=44"),
    @("1+70=71", "9+53=62"),
    @("72-38=34", "54-5=49"),
    @("74-20=54", "92-63=29"),
    @("44-36=8", "10+82=92"),
    @("19+18=37", "94-65=29"),
    @("47-26=21", "89-81=8"),
    @("64+21=85", "84-52=32"),
    @("36+63=99", "22+35=57"),
    @("87-75=12", "7+74=81"),
    @("38+19=57", "47-4=43"),
    @("97-88=9", "24+10=34"),
    @("35+16=51", "78-46=32"),
    @("67-0=67", "66+25=91"),
    @("27+17=44", "7+25=32"),
    @("15+23=38", "60-54=6"),
    @("86+11=97", "72-60=12"),
    @("38-18=20", "29+57=86"),
    @("23+53=76", "68-57=11"),
    @("2+6=8", "12+3=15"),
    @("53+16=69", "38+58=96"),
    @("52-39=13", "27+7=34"),
    @("41-30=11", "51+9=60"),
    @("43+33=76", "41+17=58"),
    @("33+28=61", "27+38=65"),
    @("36+36=72", "89-77=12"),
    @("62+17=79", "98-38=60"),
    @("95-65=30", "81-1=80"),
    @("91-85=6", "26+57=83"),
    @("48-39=9", "68-66=2"),
    @("74-12=62", "49-17=32"),
    @("60-53=7", "53-3=50"),
    @("68-43=25", "93-13=80"),
    @("35-22=13", "54-5=49"),
    @("37-13=24", "29+33=62"),
    @("40-11=29", "95-66=29"),
    @("57-2=55", "20+45=65"),
    @("89-57=32", "53-49=4"),
    @("89-21=68", "25+52=77"),
    @("25+10=35", "25-22=3"),
    @("48-11=37", "94-7=87"),
    @("61-17=44", "19+62=81"),
    @("12+81=93", "51-19=32"),
    @("12+78=90", "95-19=76"),
    @("20-20=0", "42-3=39"),
    @("40+13=53", "82-47=35"),
    @("47+33=80", "1+57=58"),
    @("62-28=34", "77+3=80"),
    @("97-86=11", "88-55=33"),
    @("63+0=63", "40-29=11"),
    @("20+38=58", "9+2=11"),
    @("34+24=58", "70-14=56"),
    @("29+28=57", "49+7=56"),
    @("93-25=68", "44+34=78"),
    @("17+41=58", "86-76=10"),
    @("7+72=79", "55-54=1"),
    @("44+49=93", "6+10=16"),
    @("79-68=11", "49+12=61"),
    @("88+7=95", "6+81=87"),
    @("6+44=50", "85-59=26"),
    @("31+3=34", "87-61=26"),
    @("22+7=29", "39+17=56"),
    @("99-72=27", "78-45=33"),
    @("80-2=78", "47-23=24"),
    @("39+27=66", "79-49=30"),
    @("0+90=90", "7+37=44"),
    @("58-17=41", "64+35=99"),
    @("44-30=14", "13+59=72"),
    @("91-85=6", "68-43=25"),
    @("56+26=82", "67-63=4"),
    @("58-47=11", "48+12=60"),
    @("5+69=74", "30+39=69"),
    @("31-22=9", "5+25=30"),
    @("1+41=42", "62-13=49"),
    @("86-56=30", "13+26=39"),
    @("31-31=0", "88-79=9"),
    @("59-37=22", "42+34=76"),
    @("0+3=3", "12+10=22"),
    @("82-17=65", "64+4=68"),
    @("41+16=57", "28+58=86"),
    @("38-8=30", "80-50=30"),
    @("8+55=63", "46-12=34"),
    @("86-75=11", "24+42=66"),
    @("12+16=28", "27-24=3"),
    @("0+32=32", "26+3=29"),
    @("47+9=56", "48-27=21"),
    @("9+3=12", "78-34=44"),
    @("69-53=16", "20+63=83"),
    @("63-26=37", "10+40=50"),
    @("76+20=96", "51+46=97"),
    @("64+2=66", "58-24=34"),
    @("5+50=55", "44-33=11"),
    @("6+25=31", "39+4=43"),
    @("49+17=66", "60-45=15"),
    @("8+23=31", "37+42=79"),
    @("75-9=66", "66+29=95"),
    @("41+23=64", "87-85=2"),
    @("27-4=23", "44-34=10"),
)

$d = $word.ActiveDocument
$t = $d.Tables(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $pair = $pairs[$idx]
        $old = $pair[0]
        $new = $pair[1]
        $cell = $t.Cell($r, $c)
        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)
        if ($current -eq $old) {
            $cell.Range.Text = $new
        } else {
            Write-Output "WARNING: cell ($r,$c) expected $old but found $current"
        }
        $idx = $idx + 1
    }
}
Write-Output "Replaced $idx cells"